$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-10) have been re-sorted by date (column D) descending,
# keeping ties in their existing relative order. Rows 4 and 5 are unchanged.
# New row <- Old row mapping (data moved between rows):
#   2 <- 9, 3 <- 10, 4 <- 4, 5 <- 5, 6 <- 8, 7 <- 3, 8 <- 2, 9 <- 7, 10 <- 6

# Capture current values for columns D, L, M, N, O, P, S for rows that move.
$srcRows = @(9, 10, 8, 3, 2, 7, 6)
$destRows = @(2, 3, 6, 7, 8, 9, 10)

$snapshot = @{}
foreach ($r in $srcRows) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

$mapping = @{ 2 = 9; 3 = 10; 6 = 8; 7 = 3; 8 = 2; 9 = 7; 10 = 6 }

foreach ($destRow in $destRows) {
    $srcRow = $mapping[$destRow]
    $vals = $snapshot[$srcRow]
    $ws.Cells.Item($destRow, 4).Value2 = $vals.D
    $ws.Cells.Item($destRow, 12).Value2 = $vals.L
    $ws.Cells.Item($destRow, 13).Value2 = $vals.M
    $ws.Cells.Item($destRow, 14).Value2 = $vals.N
    $ws.Cells.Item($destRow, 15).Value2 = $vals.O
    $ws.Cells.Item($destRow, 16).Value2 = $vals.P
    $ws.Cells.Item($destRow, 19).Value2 = $vals.S
}
